$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tareas diarias")

# ---- Row 8: "Completar el proceso de compra" -> estado pasa de Incompleto a Trabajando ----
# La fila 17 ya tiene el formato "Trabajando" (estilos 6/3/3/8) que necesitamos copiar.
$ws.Range("A17:C17").Copy() | Out-Null
$ws.Range("A8:C8").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(8).AutoFit() | Out-Null
$ws.Range("D8").Value = "Al finalizar mostrar una factura, revisar si le llega al administrador"

# ---- Row 10: "Resolver simbolos raros..." -> estado pasa de Incompleto a Completo ----
# La fila 7 ya tiene el formato "Completo" (estilos 5/2/2/7) que necesitamos copiar.
$ws.Range("A7:D7").Copy() | Out-Null
$ws.Range("A10:D10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value = "Completo"
$ws.Range("D10").Value = "fixed ñ + vocales con tildes"

# ---- Row 25: se completa con Responsable/Estado y se ajusta el formato ----
$ws.Range("A23").Copy() | Out-Null
$ws.Range("A25").PasteSpecial(-4122) | Out-Null
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B25:C25").PasteSpecial(-4122) | Out-Null
$ws.Range("B25").Value = "Fanky"
$ws.Range("C25").Value = "Completo"
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null

# ---- sheetView: se quita el topLeftCell y la seleccion pasa a D18 ----
$ws.Activate() | Out-Null
$ws.Range("D18").Select() | Out-Null
